# Apply ObjTables schema/date bump + newly-populated "verbose name" column
# (camelCase attribute names) across the workbook, per commit:
#   renaming EmpiricalFormulaAttribute->ChemicalFormulaAttribute; enabling
#   serializing *-to-many relationships to CSV/TSV; adding
#   ReactionEquationAttribute and ListAttribute; adding metabolomics example.
#
# Concretely, the template.xlsx header stamps move from
#   objTablesVersion='0.0.9' date='2020-04-27 01:05:01'
# to
#   objTablesVersion='1.0.0' date='2020-05-29 00:18:53'
# and the previously-blank "!Verbose name" column on the "!!_Schema" sheet
# is now populated with the CamelCase verbose name for every attribute, and
# the Slug(...) regex attribute's definition is switched to a raw string
# literal (r'...') with single backslashes instead of escaped double ones.

$wb = $excel.ActiveWorkbook

$oldVersion = "0.0.9"
$newVersion = "1.0.0"
$oldDate = "2020-04-27 01:05:01"
$newDate = "2020-05-29 00:18:53"

# ---------------------------------------------------------------------
# Sheet "!!_Table of contents"
# ---------------------------------------------------------------------
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
$wsToc.Range("A1").Value = "!!!ObjTables objTablesVersion='$newVersion' date='$newDate'"
$wsToc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='$newDate' objTablesVersion='$newVersion'"

# ---------------------------------------------------------------------
# Sheet "!!_Schema"
# ---------------------------------------------------------------------
$wsSchema = $wb.Worksheets.Item("!!_Schema")
$wsSchema.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='$newDate' objTablesVersion='$newVersion'"

# id Slug(...) attribute definition: switch to raw-string regex literal
$slugDef = 'Slug(r''^(?!(^|\b)(\d+(\.\d*)?(\b|$))|(\.\d+$)|(0[x][0-9a-f]+(\b|$))|([0-9]+e[0-9]+(\b|$)))[a-z0-9_]+$'', flags=2, primary=True, unique=True)'
$wsSchema.Range("D4").Value = $slugDef
$wsSchema.Range("D10").Value = $slugDef
$wsSchema.Range("D15").Value = $slugDef

# Populate the "!Verbose name" column (E) for every attribute row
$wsSchema.Range("E4").Value = "Id"
$wsSchema.Range("E5").Value = "Identifiers"
$wsSchema.Range("E6").Value = "IsConstant"
$wsSchema.Range("E7").Value = "Model"
$wsSchema.Range("E8").Value = "Name"

$wsSchema.Range("E10").Value = "Id"
$wsSchema.Range("E11").Value = "Name"

$wsSchema.Range("E13").Value = "Equation"
$wsSchema.Range("E14").Value = "Gene"
$wsSchema.Range("E15").Value = "Id"
$wsSchema.Range("E16").Value = "Identifiers"
$wsSchema.Range("E17").Value = "IsReversible"
$wsSchema.Range("E18").Value = "Model"
$wsSchema.Range("E19").Value = "Name"

# ---------------------------------------------------------------------
# Sheet "!!Compound"
# ---------------------------------------------------------------------
$wsCompound = $wb.Worksheets.Item("!!Compound")
$wsCompound.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Compound' name='Compound' description='Compound' date='$newDate' objTablesVersion='$newVersion'"

# ---------------------------------------------------------------------
# Sheet "!!Model"
# ---------------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("!!Model")
$wsModel.Range("A1").Value = "!!ObjTables type='Data' tableFormat='column' class='Model' name='Model' description='Model' date='$newDate' objTablesVersion='$newVersion'"

# ---------------------------------------------------------------------
# Sheet "!!Reaction"
# ---------------------------------------------------------------------
$wsReaction = $wb.Worksheets.Item("!!Reaction")
$wsReaction.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Reaction' name='Reaction' description='Reaction' date='$newDate' objTablesVersion='$newVersion'"
